$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.531.55"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "3.102.54"
$ws.Range("E3").Value = "  +2.67%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "385.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0855"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").Value = "3.596.63"
$ws.Range("E13").Value = "  +2.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "3.098.81"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.998"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.72%  "
$ws.Range("D19").Value = "51.583.20"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.43%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -2.92%  "
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0466"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.296"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.18%  "
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "129.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.116"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("E47").Value = "  +4.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").Value = "2.064.89"
$ws.Range("E49").Value = "  +1.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0330"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.901"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.39%  "
